# Update "想去人数" (F column) counts across the four sheets of the
# 广州-漫展信息 workbook to match the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value  = 304
$ws.Cells.Item(4, 6).Value  = 1238
$ws.Cells.Item(5, 6).Value  = 352
$ws.Cells.Item(6, 6).Value  = 317
$ws.Cells.Item(7, 6).Value  = 3821
$ws.Cells.Item(9, 6).Value  = 744
$ws.Cells.Item(10, 6).Value = 1887
$ws.Cells.Item(11, 6).Value = 327
$ws.Cells.Item(12, 6).Value = 215
$ws.Cells.Item(13, 6).Value = 729
$ws.Cells.Item(14, 6).Value = 149
$ws.Cells.Item(16, 6).Value = 2081
$ws.Cells.Item(18, 6).Value = 4
$ws.Cells.Item(21, 6).Value = 221

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value  = 38
$ws.Cells.Item(7, 6).Value  = 31
$ws.Cells.Item(12, 6).Value = 80
$ws.Cells.Item(22, 6).Value = 48
$ws.Cells.Item(23, 6).Value = 55

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 6390
$ws.Cells.Item(4, 6).Value = 2068

# Sheet "全部类型" (All Types) - combined list of every event above
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value  = 6390
$ws.Cells.Item(4, 6).Value  = 2068
$ws.Cells.Item(9, 6).Value  = 38
$ws.Cells.Item(10, 6).Value = 38
$ws.Cells.Item(12, 6).Value = 304
$ws.Cells.Item(13, 6).Value = 1238
$ws.Cells.Item(14, 6).Value = 353
$ws.Cells.Item(16, 6).Value = 31
$ws.Cells.Item(18, 6).Value = 317
$ws.Cells.Item(19, 6).Value = 3821
$ws.Cells.Item(24, 6).Value = 80
$ws.Cells.Item(25, 6).Value = 744
$ws.Cells.Item(26, 6).Value = 1887
$ws.Cells.Item(27, 6).Value = 327
$ws.Cells.Item(29, 6).Value = 215
$ws.Cells.Item(30, 6).Value = 729
$ws.Cells.Item(31, 6).Value = 149
$ws.Cells.Item(34, 6).Value = 2081
$ws.Cells.Item(38, 6).Value = 4
$ws.Cells.Item(41, 6).Value = 221
$ws.Cells.Item(48, 6).Value = 48
$ws.Cells.Item(49, 6).Value = 55
